$wb = $excel.ActiveWorkbook

# This report records a new handoff event for the
# "cc0121d2-f0d0-4fec-8829-8fd6ee143b4b" file (row 3 on each sheet):
#  - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#  - The per-locale "Latest Handoff Datetime" is stamped with the new handoff time

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-10 07:09:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-10 07:09:48"
